$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down one row
$ws.Rows.Item(1).Insert()

# Populate the new header row (order matches shared-string insertion order:
# Year, Change, Market cap)
$ws.Range("A1").Value = "Year"
$ws.Range("C1").Value = "Change"
$ws.Range("B1").Value = "Market cap"

# Update selection to match the post-edit state
$ws.Range("E6").Select()
